# Saldo.xlsx ("Export" sheet) maintenance edit:
#   - Fold a few closed/zeroed account balances into the account they were
#     reconciled against (two running-balance cells get new totals).
#   - Remove the now-redundant rows for the folded-in accounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Account -> new balance for the two rows that absorb the removed accounts.
$balanceUpdates = @{
    "004267119" = 23179.66   # Ana
    "008115265" = 8345.02    # Elaine
}

foreach ($acct in $balanceUpdates.Keys) {
    $cell = $ws.Columns(1).Find($acct)
    if ($cell -eq $null) {
        throw "Account $acct not found in column A"
    }
    $ws.Cells.Item($cell.Row, 3).Value = $balanceUpdates[$acct]
}

# Accounts whose rows are being removed outright (folded into the totals above).
$accountsToRemove = @("008030888", "004474776", "004461526", "004204344", "004383268", "004361159")

$rowsToRemove = @()
foreach ($acct in $accountsToRemove) {
    $cell = $ws.Columns(1).Find($acct)
    if ($cell -eq $null) {
        throw "Account $acct not found in column A"
    }
    $rowsToRemove += $cell.Row
}

# Delete bottom-up so earlier row numbers stay valid as rows shift up.
$rowsToRemove = $rowsToRemove | Sort-Object -Descending
foreach ($r in $rowsToRemove) {
    $ws.Rows($r).Delete()
}
